$wb = $excel.ActiveWorkbook

# Copy the Portugal sheet (last sheet) to create the new Slovakia sheet, placed after it
$portugal = $wb.Worksheets.Item("Portugal")
$portugal.Copy($null, $portugal)

# The newly copied sheet becomes the active sheet, placed right after Portugal
$newSheet = $wb.ActiveSheet
$newSheet.Name = "Slovakia"

# Update values
$newSheet.Range("B2").Value = "Slovakia Market"
$newSheet.Range("B4").Value = "NGC-2930/T3174"
$newSheet.Range("B4").Style = "Normal"

# Remove the custom row heights inherited from Portugal (rows 3-5) so the new
# sheet reverts to the sheet default row height
$newSheet.Rows("3:5").AutoFit()

# Restore Portugal's selection state (mirrors Excel's "Move or Copy" side-effect
# of leaving the whole sheet selected on the source sheet)
$portugal.Cells.Select()

# Re-activate the new Slovakia sheet/tab and set its own selection
$newSheet.Activate()
$newSheet.Range("B4").Select()
